$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(2)
$tf = $shp.TextFrame
$tr = $tf.TextRange

# Paragraph 3: "Current " + "Job Openings:" -> "Current Job Openings:"
$para3 = $tr.Paragraphs(3, 1)
$run1 = $para3.Runs(1, 1)
$run2 = $para3.Runs(2, 1)
$run1.Text = "Current Job Openings:"
$run2.Text = ""

# Paragraph 6: "...WebKit" + " " + "Development, Quality and automation"
#   -> merge the " " run and the "Development, Quality and automation" run
$para6 = $tr.Paragraphs(6, 1)
$run3 = $para6.Runs(3, 1)
$run4 = $para6.Runs(4, 1)
$run3.Text = " Development, Quality and automation"
$run4.Text = ""
